# Commit "#5: fund, bonds, otherbonds, antique done"
#
# Sheet 5 (基金受益憑證 / "funds") gets turned into a proper scraped table,
# matching the shape the other sheets (e.g. 股票/stock on sheet 4) already
# have:
#   - Row 1 used to duplicate the first data row's values; it becomes a
#     real header row (name/owner/dealer/quantity/face_value/currency/total
#     plus the scraper metadata headers).
#   - Columns I:O are added to every data row with the scraper metadata:
#       property_category, category, date, legislator_name, legislator_id,
#       source_file, index

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基金受益憑證")

# --- Row 1: turn the duplicated data row into a header row ---------------
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "owner"
$ws.Cells.Item(1, 4).Value = "dealer"
$ws.Cells.Item(1, 5).Value = "quantity"
$ws.Cells.Item(1, 6).Value = "face_value"
$ws.Cells.Item(1, 7).Value = "currency"
$ws.Cells.Item(1, 8).Value = "total"
$ws.Cells.Item(1, 9).Value = "property_category"
$ws.Cells.Item(1, 10).Value = "category"
$ws.Cells.Item(1, 11).Value = "date"
$ws.Cells.Item(1, 12).Value = "legislator_name"
$ws.Cells.Item(1, 13).Value = "legislator_id"
$ws.Cells.Item(1, 14).Value = "source_file"
$ws.Cells.Item(1, 15).Value = "index"

# --- Rows 2 & 3: append the scraper metadata columns (I:O) ---------------
# The "date" column (K) holds a literal "2012-04-27" string, not a real
# date - force the cell to text first so COM doesn't coerce it to a date
# serial number.
$ws.Cells.Item(2, 11).NumberFormat = "@"
$ws.Cells.Item(3, 11).NumberFormat = "@"

# Row 2 (record index 95)
$ws.Cells.Item(2, 9).Value = "fund"
$ws.Cells.Item(2, 10).Value = "normal"
$ws.Cells.Item(2, 11).Value = "2012-04-27"
$ws.Cells.Item(2, 12).Value = "江惠貞"
$ws.Cells.Item(2, 13).Value = 1732
$ws.Cells.Item(2, 14).Value = "tmpf6f41"
$ws.Cells.Item(2, 15).Value = 95

# Row 3 (record index 96)
$ws.Cells.Item(3, 9).Value = "fund"
$ws.Cells.Item(3, 10).Value = "normal"
$ws.Cells.Item(3, 11).Value = "2012-04-27"
$ws.Cells.Item(3, 12).Value = "江惠貞"
$ws.Cells.Item(3, 13).Value = 1732
$ws.Cells.Item(3, 14).Value = "tmpf6f41"
$ws.Cells.Item(3, 15).Value = 96

Write-Output "sheet5 updated"
